$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text formatting so numeric-looking price strings
# (e.g. "1.00", "0.628") are not silently coerced into numbers by Excel,
# matching the source data which stores prices as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.221.71'
$ws.Range("E2").Value = '  -4.24%  '

$ws.Range("D3").Value = '3.708.83'
$ws.Range("E3").Value = '  -4.62%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").Value = '597.51'
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").Value = '183.06'
$ws.Range("E6").Value = '  +8.22%  '

$ws.Range("D7").Value = '3.702.24'
$ws.Range("E7").Value = '  -4.51%  '

$ws.Range("D8").Value = '0.628'
$ws.Range("E8").Value = '  -6.71%  '

$ws.Range("D9").Value = '0.996'
$ws.Range("E9").Value = '  -0.43%  '

$ws.Range("E10").Value = '  -5.18%  '

$ws.Range("E11").Value = '  -7.64%  '

$ws.Range("D12").Value = '56.13'
$ws.Range("E12").Value = '  +4.20%  '

$ws.Range("E13").Value = '  -10.45%  '

$ws.Range("E14").Value = '  -8.73%  '

$ws.Range("D15").Value = '4.274.51'
$ws.Range("E15").Value = '  -4.79%  '

$ws.Range("D16").Value = '3.699.47'
$ws.Range("E16").Value = '  -4.00%  '

$ws.Range("D17").Value = '19.44'
$ws.Range("E17").Value = '  -6.85%  '

$ws.Range("E18").Value = '  -2.30%  '

$ws.Range("D19").Value = '12.85'
$ws.Range("E19").Value = '  -7.63%  '

$ws.Range("D20").Value = '1.13'
$ws.Range("E20").Value = '  -7.04%  '

$ws.Range("D21").Value = '68.098.27'
$ws.Range("E21").Value = '  -3.88%  '

$ws.Range("D22").Value = '409.87'
$ws.Range("E22").Value = '  -6.41%  '

$ws.Range("D23").Value = '4.57'
$ws.Range("E23").Value = '  -3.50%  '

$ws.Range("D24").Value = '88.69'
$ws.Range("E24").Value = '  -6.43%  '

$ws.Range("D25").Value = '3.02'
$ws.Range("E25").Value = '  -8.49%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '12.81'
$ws.Range("E26").Value = '  -7.79%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '11.02'
$ws.Range("E27").Value = '  -2.28%  '

$ws.Range("D28").Value = '3.86'
$ws.Range("E28").Value = '  -6.17%  '

$ws.Range("D29").Value = '6.07'
$ws.Range("E29").Value = '  +2.36%  '

$ws.Range("D30").Value = '9.51'
$ws.Range("E30").Value = '  -8.43%  '

$ws.Range("D31").Value = '32.81'
$ws.Range("E31").Value = '  -6.96%  '

$ws.Range("D32").Value = '7.28'
$ws.Range("E32").Value = '  -10.60%  '

$ws.Range("D33").Value = '12.52'
$ws.Range("E33").Value = '  -8.07%  '

$ws.Range("E34").Value = '  -6.81%  '

$ws.Range("D35").Value = '43.71'
$ws.Range("E35").Value = '  -11.16%  '

$ws.Range("D36").Value = '64.18'
$ws.Range("E36").Value = '  -9.10%  '

$ws.Range("D37").Value = '592.59'
$ws.Range("E37").Value = '  -6.33%  '

$ws.Range("D38").Value = '0.0₃0883'
$ws.Range("E38").Value = '  -10.93%  '

$ws.Range("E39").Value = '  -5.81%  '

$ws.Range("E40").Value = '  +0.36%  '

$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  -5.54%  '

$ws.Range("D43").Value = '2.76'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").Value = '3.01'
$ws.Range("E44").Value = '  -8.07%  '

$ws.Range("D45").Value = '0.0436'
$ws.Range("E45").Value = '  -7.66%  '

$ws.Range("D46").Value = '2.90'
$ws.Range("E46").Value = '  -14.42%  '

$ws.Range("D47").Value = '9.22'
$ws.Range("E47").Value = '  -9.10%  '

$ws.Range("D48").Value = '2.72'
$ws.Range("E48").Value = '  -3.96%  '

$ws.Range("E49").Value = '  -6.98%  '

$ws.Range("D50").Value = '2.771.28'
$ws.Range("E50").Value = '  -2.07%  '

$ws.Range("E51").Value = '  -5.08%  '
